$d = $word.ActiveDocument

# Locate the run that holds "dev分支。" (the text right after the bookmark)
# and the end of the preceding run ("...创建了一个").
$findRng = $d.Content.Duplicate
$findRng.Find.ClearFormatting()
$findRng.Find.Execute("dev分支。", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($findRng.Find.Found) {
    $runStart = $findRng.Start
    $runEnd = $findRng.End

    # Remove the whole "dev分支。" run entirely (element and all).
    $oldRun = $d.Range($runStart, $runEnd)
    $oldRun.Delete()

    # Append the combined new text right after the prior run
    # ("...创建了一个"), so it merges into that run exactly like the
    # diff shows (single run, no new run node created).
    $insertPoint = $d.Range($runStart, $runStart)
    $insertPoint.InsertAfter("dev分支。使用git创建分支简单又快速。")
}
